# Auto-generated Excel COM-interop script to apply scheduled market-data update
# to the Brynhildr_Profits workbook across all 8 profession sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Updates currentAveragePrice / currentAveragePriceNQ /
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H,I,J,K,L,M,N) with refreshed market board values.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1009.5833
$ws.Range("J18").Value = 899.6
$ws.Range("L18").Value = 899.6
$ws.Range("N18").Value = -1467.6
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H51").Value = 5323.4116
$ws.Range("I51").Value = 4999.9375
$ws.Range("K51").Value = 4999.9375
$ws.Range("M51").Value = -4515.9375
$ws.Range("H103").Value = 599
$ws.Range("I103").Value = 599
$ws.Range("K103").Value = 1797
$ws.Range("M103").Value = -1211
$ws.Range("H113").Value = 3583.25
$ws.Range("I113").Value = 2999.375
$ws.Range("K113").Value = 2999.375
$ws.Range("M113").Value = 254.625
$ws.Range("H132").Value = 11409.5
$ws.Range("I132").Value = 13663.944
$ws.Range("K132").Value = 40991.83199999999
$ws.Range("M132").Value = -38461.83199999999
$ws.Range("H135").Value = 3971.037
$ws.Range("I135").Value = 1011.5263
$ws.Range("K135").Value = 9103.736699999999
$ws.Range("M135").Value = -6568.736699999999
$ws.Range("H137").Value = 15388836
$ws.Range("I137").Value = 33334840
$ws.Range("K137").Value = 100004520
$ws.Range("M137").Value = -100001970
$ws.Range("H138").Value = 3496.0278
$ws.Range("I138").Value = 2523.111
$ws.Range("J138").Value = 3820.3333
$ws.Range("K138").Value = 7569.333
$ws.Range("L138").Value = 11460.9999
$ws.Range("M138").Value = -2429.333
$ws.Range("N138").Value = -21740.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1187754.8
$ws.Range("I74").Value = 1464427.6
$ws.Range("K74").Value = 1464427.6
$ws.Range("M74").Value = -1463553.6
$ws.Range("H77").Value = 1187754.8
$ws.Range("I77").Value = 1464427.6
$ws.Range("K77").Value = 7322138
$ws.Range("M77").Value = -7317770
$ws.Range("H88").Value = 3255.8
$ws.Range("I88").Value = 2173.5
$ws.Range("J88").Value = 4338.1
$ws.Range("K88").Value = 2173.5
$ws.Range("L88").Value = 4338.1
$ws.Range("M88").Value = -1767.5
$ws.Range("N88").Value = -5150.1
$ws.Range("H91").Value = 3255.8
$ws.Range("I91").Value = 2173.5
$ws.Range("J91").Value = 4338.1
$ws.Range("K91").Value = 2173.5
$ws.Range("L91").Value = 4338.1
$ws.Range("M91").Value = -769.5
$ws.Range("N91").Value = -7146.1
$ws.Range("H109").Value = 69501.5
$ws.Range("J109").Value = 69501.5
$ws.Range("L109").Value = 69501.5
$ws.Range("N109").Value = -72275.5
$ws.Range("H132").Value = 679535.2
$ws.Range("I132").Value = 785153.25
$ws.Range("J132").Value = 3579.6
$ws.Range("K132").Value = 2355459.75
$ws.Range("L132").Value = 10738.8
$ws.Range("M132").Value = -2352929.75
$ws.Range("N132").Value = -15798.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3037092.8
$ws.Range("I134").Value = 4863.564
$ws.Range("K134").Value = 14590.692
$ws.Range("M134").Value = -12055.692

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 78808.84
$ws.Range("I16").Value = 112501.78
$ws.Range("J16").Value = 2999.75
$ws.Range("K16").Value = 112501.78
$ws.Range("L16").Value = 2999.75
$ws.Range("M16").Value = -112214.78
$ws.Range("N16").Value = -3573.75
$ws.Range("H105").Value = 5687.087
$ws.Range("I105").Value = 6612.5293
$ws.Range("J105").Value = 3065
$ws.Range("K105").Value = 6612.5293
$ws.Range("L105").Value = 3065
$ws.Range("M105").Value = -4865.5293
$ws.Range("N105").Value = -6559
$ws.Range("H113").Value = 78808.84
$ws.Range("I113").Value = 112501.78
$ws.Range("J113").Value = 2999.75
$ws.Range("K113").Value = 112501.78
$ws.Range("L113").Value = 2999.75
$ws.Range("M113").Value = -110331.78
$ws.Range("N113").Value = -7339.75
$ws.Range("H122").Value = 27026
$ws.Range("I122").Value = 4030.3333
$ws.Range("K122").Value = 12090.9999
$ws.Range("M122").Value = -9640.999899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2565399.8
$ws.Range("I5").Value = 4465058
$ws.Range("J5").Value = 1552248.8
$ws.Range("K5").Value = 13395174
$ws.Range("L5").Value = 4656746.4
$ws.Range("M5").Value = -13395062
$ws.Range("N5").Value = -4656970.4
$ws.Range("H34").Value = 2799.3333
$ws.Range("J34").Value = 3124.25
$ws.Range("L34").Value = 9372.75
$ws.Range("N34").Value = -9540.75
$ws.Range("H68").Value = 6426.269
$ws.Range("J68").Value = 8888.777
$ws.Range("L68").Value = 26666.331
$ws.Range("N68").Value = -28288.331
$ws.Range("H71").Value = 6426.269
$ws.Range("J71").Value = 8888.777
$ws.Range("L71").Value = 79998.993
$ws.Range("N71").Value = -88110.993
$ws.Range("H80").Value = 11435.728
$ws.Range("J80").Value = 13499.125
$ws.Range("L80").Value = 40497.375
$ws.Range("N80").Value = -42369.375
$ws.Range("H83").Value = 11435.728
$ws.Range("J83").Value = 13499.125
$ws.Range("L83").Value = 121492.125
$ws.Range("N83").Value = -130852.125
$ws.Range("H114").Value = 8173.1333
$ws.Range("J114").Value = 10405.565
$ws.Range("L114").Value = 31216.695
$ws.Range("N114").Value = -37724.695
$ws.Range("H120").Value = 26874.875
$ws.Range("I120").Value = 24999
$ws.Range("K120").Value = 74997
$ws.Range("M120").Value = -70159
$ws.Range("H121").Value = 3580896.2
$ws.Range("I121").Value = 558
$ws.Range("J121").Value = 4557352
$ws.Range("K121").Value = 1674
$ws.Range("L121").Value = 13672056
$ws.Range("M121").Value = -364
$ws.Range("N121").Value = -13674676
$ws.Range("H135").Value = 2565399.8
$ws.Range("I135").Value = 4465058
$ws.Range("J135").Value = 1552248.8
$ws.Range("K135").Value = 40185522
$ws.Range("L135").Value = 13970239.2
$ws.Range("M135").Value = -40182987
$ws.Range("N135").Value = -13975309.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 15540.083
$ws.Range("J43").Value = 38332.332
$ws.Range("L43").Value = 38332.332
$ws.Range("N43").Value = -38634.332
$ws.Range("H97").Value = 696.64703
$ws.Range("I97").Value = 886.3
$ws.Range("J97").Value = 425.7143
$ws.Range("K97").Value = 886.3
$ws.Range("L97").Value = 425.7143
$ws.Range("M97").Value = -390.3
$ws.Range("N97").Value = -1417.7143
$ws.Range("H122").Value = 6199.467
$ws.Range("I122").Value = 9057.799999999999
$ws.Range("J122").Value = 4770.3
$ws.Range("K122").Value = 27173.4
$ws.Range("L122").Value = 14310.9
$ws.Range("M122").Value = -24723.4
$ws.Range("N122").Value = -19210.9
$ws.Range("H132").Value = 12540.763
$ws.Range("I132").Value = 11366.226
$ws.Range("J132").Value = 17742.285
$ws.Range("K132").Value = 34098.678
$ws.Range("L132").Value = 53226.855
$ws.Range("M132").Value = -31568.678
$ws.Range("N132").Value = -58286.855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3341141.5
$ws.Range("I132").Value = 4496404.5
$ws.Range("J132").Value = 3715.889
$ws.Range("K132").Value = 13489213.5
$ws.Range("L132").Value = 11147.667
$ws.Range("M132").Value = -13486683.5
$ws.Range("N132").Value = -16207.667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 4000
$ws.Range("J7").Value = 4000
$ws.Range("L7").Value = 4000
$ws.Range("N7").Value = -4226
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 119411.6
$ws.Range("I122").Value = 11380
$ws.Range("J122").Value = 227443.2
$ws.Range("K122").Value = 34140
$ws.Range("L122").Value = 682329.6000000001
$ws.Range("M122").Value = -31690
$ws.Range("N122").Value = -687229.6000000001
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()
